# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Replaces the 38 "Periodo Mora" data rows (rows 16-53, columns E:G) with the
# new, chronologically-ascending set of periods (1702 .. 2003) and their
# updated "Valor Mora" (F) / "Salario Basico" (G) amounts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @(
    "1702", "1703", "1704", "1705", "1706", "1707", "1708", "1709", "1710", "1711", "1712",
    "1801", "1802", "1803", "1804", "1805", "1806", "1807", "1808", "1809", "1810", "1811", "1812",
    "1901", "1902", "1903", "1904", "1905", "1906", "1907", "1908", "1909", "1910", "1911", "1912",
    "2001", "2002", "2003"
)

$firstRow = 16

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $firstRow + $i

    $ws.Cells.Item($row, 5).Value = $periods[$i]

    if ($row -le 34) {
        $ws.Cells.Item($row, 6).Value = 27578
    } else {
        $ws.Cells.Item($row, 6).Value = 31249
    }

    $ws.Cells.Item($row, 7).Value = 781242
}
